$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 737.56525
$ws.Range("I41").Value = 1486.6666
$ws.Range("K41").Value = 1486.6666
$ws.Range("M41").Value = -1046.6666
$ws.Range("H96").Value = 1949953.2
$ws.Range("I96").Value = 302.66666
$ws.Range("J96").Value = 3509673.8
$ws.Range("K96").Value = 907.9999799999999
$ws.Range("L96").Value = 10529021.4
$ws.Range("M96").Value = 465.0000200000001
$ws.Range("N96").Value = -10531767.4
$ws.Range("H125").Value = 3833.3333
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 4750
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 42750
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -47670
$ws.Range("H133").Value = 48249.9
$ws.Range("J133").Value = 48249.9
$ws.Range("L133").Value = 48249.9
$ws.Range("N133").Value = -58369.9
$ws.Range("H134").Value = 48099.6
$ws.Range("J134").Value = 48099.6
$ws.Range("L134").Value = 48099.6
$ws.Range("N134").Value = -58239.6
$ws.Range("H136").Value = 47599.8
$ws.Range("J136").Value = 47599.8
$ws.Range("L136").Value = 47599.8
$ws.Range("N136").Value = -57799.8
$ws.Range("H139").Value = 33811.8
$ws.Range("J139").Value = 33811.8
$ws.Range("L139").Value = 33811.8
$ws.Range("N139").Value = -44091.8
$ws.Range("H140").Value = 35969.75
$ws.Range("J140").Value = 35969.75
$ws.Range("L140").Value = 35969.75
$ws.Range("N140").Value = -46329.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10727.738
$ws.Range("I32").Value = 4022.9512
$ws.Range("J32").Value = 22181.75
$ws.Range("K32").Value = 4022.9512
$ws.Range("L32").Value = 22181.75
$ws.Range("M32").Value = -3735.9512
$ws.Range("N32").Value = -22755.75
$ws.Range("H61").Value = 2567.6667
$ws.Range("I61").Value = 1687
$ws.Range("J61").Value = 5650
$ws.Range("K61").Value = 1687
$ws.Range("L61").Value = 5650
$ws.Range("M61").Value = -1475
$ws.Range("N61").Value = -6074
$ws.Range("H74").Value = 20558162
$ws.Range("I74").Value = 18002528
$ws.Range("J74").Value = 33336336
$ws.Range("K74").Value = 18002528
$ws.Range("L74").Value = 33336336
$ws.Range("M74").Value = -18001654
$ws.Range("N74").Value = -33338084
$ws.Range("H77").Value = 20558162
$ws.Range("I77").Value = 18002528
$ws.Range("J77").Value = 33336336
$ws.Range("K77").Value = 90012640
$ws.Range("L77").Value = 166681680
$ws.Range("M77").Value = -90008272
$ws.Range("N77").Value = -166690416
$ws.Range("H86").Value = 30000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 30000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 30000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -32372
$ws.Range("H89").Value = 30000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 30000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 90000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -101856
$ws.Range("H136").Value = 2567.6667
$ws.Range("I136").Value = 1687
$ws.Range("J136").Value = 5650
$ws.Range("K136").Value = 5061
$ws.Range("L136").Value = 16950
$ws.Range("M136").Value = -2511
$ws.Range("N136").Value = -22050

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3335.889
$ws.Range("I107").Value = 3011.6667
$ws.Range("J107").Value = 3984.3333
$ws.Range("K107").Value = 3011.6667
$ws.Range("L107").Value = 3984.3333
$ws.Range("M107").Value = -1091.6667
$ws.Range("N107").Value = -7824.3333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 297.27274
$ws.Range("I22").Value = 247.5
$ws.Range("K22").Value = 247.5
$ws.Range("M22").Value = 102.5
$ws.Range("H31").Value = 8002467.5
$ws.Range("I31").Value = 16669333
$ws.Range("J31").Value = 5265562.5
$ws.Range("K31").Value = 16669333
$ws.Range("L31").Value = 5265562.5
$ws.Range("M31").Value = -16669038
$ws.Range("N31").Value = -5266152.5
$ws.Range("H34").Value = 8002467.5
$ws.Range("I34").Value = 16669333
$ws.Range("J34").Value = 5265562.5
$ws.Range("K34").Value = 16669333
$ws.Range("L34").Value = 5265562.5
$ws.Range("M34").Value = -16669131
$ws.Range("N34").Value = -5265966.5
$ws.Range("H68").Value = 17999.666
$ws.Range("J68").Value = 17999.666
$ws.Range("L68").Value = 17999.666
$ws.Range("N68").Value = -19497.666
$ws.Range("H71").Value = 17999.666
$ws.Range("J71").Value = 17999.666
$ws.Range("L71").Value = 53998.99800000001
$ws.Range("N71").Value = -61486.99800000001
$ws.Range("H122").Value = 2565582.5
$ws.Range("I122").Value = 4762844
$ws.Range("J122").Value = 2110.5
$ws.Range("K122").Value = 14288532
$ws.Range("L122").Value = 6331.5
$ws.Range("M122").Value = -14286082
$ws.Range("N122").Value = -11231.5
$ws.Range("H132").Value = 2042.6285
$ws.Range("I132").Value = 1045.9166
$ws.Range("J132").Value = 4217.273
$ws.Range("K132").Value = 3137.7498
$ws.Range("L132").Value = 12651.819
$ws.Range("M132").Value = -607.7498000000001
$ws.Range("N132").Value = -17711.819

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 750.26
$ws.Range("I68").Value = 391.31033
$ws.Range("J68").Value = 1245.9524
$ws.Range("K68").Value = 1173.93099
$ws.Range("L68").Value = 3737.857199999999
$ws.Range("M68").Value = -362.9309900000001
$ws.Range("N68").Value = -5359.857199999999
$ws.Range("H71").Value = 750.26
$ws.Range("I71").Value = 391.31033
$ws.Range("J71").Value = 1245.9524
$ws.Range("K71").Value = 3521.79297
$ws.Range("L71").Value = 11213.5716
$ws.Range("M71").Value = 534.2070299999996
$ws.Range("N71").Value = -19325.5716
$ws.Range("H87").Value = 4941.3794
$ws.Range("I87").Value = 766.6667
$ws.Range("J87").Value = 5423.077
$ws.Range("K87").Value = 2300.0001
$ws.Range("L87").Value = 16269.231
$ws.Range("M87").Value = -1052.0001
$ws.Range("N87").Value = -18765.231
$ws.Range("H90").Value = 4941.3794
$ws.Range("I90").Value = 766.6667
$ws.Range("J90").Value = 5423.077
$ws.Range("K90").Value = 6900.0003
$ws.Range("L90").Value = 48807.693
$ws.Range("M90").Value = -660.0002999999997
$ws.Range("N90").Value = -61287.693
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1060
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1383.6666
$ws.Range("I132").Value = 434
$ws.Range("K132").Value = 3906
$ws.Range("M132").Value = -1376

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2177.111
$ws.Range("I122").Value = 2416.2856
$ws.Range("J122").Value = 1340
$ws.Range("K122").Value = 7248.8568
$ws.Range("L122").Value = 4020
$ws.Range("M122").Value = -4798.8568
$ws.Range("N122").Value = -8920

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4019.6875
$ws.Range("I7").Value = 3954.5454
$ws.Range("J7").Value = 4163
$ws.Range("K7").Value = 3954.5454
$ws.Range("L7").Value = 4163
$ws.Range("M7").Value = -3842.5454
$ws.Range("N7").Value = -4387
$ws.Range("H22").Value = 1754870.9
$ws.Range("I22").Value = 4166984
$ws.Range("J22").Value = 606.8182
$ws.Range("K22").Value = 4166984
$ws.Range("L22").Value = 606.8182
$ws.Range("M22").Value = -4166689
$ws.Range("N22").Value = -1196.8182
$ws.Range("H27").Value = 1754870.9
$ws.Range("I27").Value = 4166984
$ws.Range("J27").Value = 606.8182
$ws.Range("K27").Value = 4166984
$ws.Range("L27").Value = 606.8182
$ws.Range("M27").Value = -4166877
$ws.Range("N27").Value = -820.8182
$ws.Range("H126").Value = 4019.6875
$ws.Range("I126").Value = 3954.5454
$ws.Range("J126").Value = 4163
$ws.Range("K126").Value = 11863.6362
$ws.Range("L126").Value = 12489
$ws.Range("M126").Value = -9393.636200000001
$ws.Range("N126").Value = -17429
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 47550
$ws.Range("J131").Value = 47550
$ws.Range("L131").Value = 47550
$ws.Range("N131").Value = -57630

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 62501360
$ws.Range("I113").Value = 45456160
$ws.Range("K113").Value = 136368480
$ws.Range("M113").Value = -136366310
$ws.Range("H122").Value = 55556990
$ws.Range("I122").Value = 90909970
$ws.Range("J122").Value = 2314.2856
$ws.Range("K122").Value = 272729910
$ws.Range("L122").Value = 6942.8568
$ws.Range("M122").Value = -272727460
$ws.Range("N122").Value = -11842.8568
$ws.Range("H130").Value = 60976.332
$ws.Range("J130").Value = 60976.332
$ws.Range("L130").Value = 60976.332
$ws.Range("N130").Value = -71016.33199999999
$ws.Range("H131").Value = 35000
$ws.Range("J131").Value = 35000
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080
$ws.Range("H132").Value = 1774.3793
$ws.Range("I132").Value = 1308.1052
$ws.Range("J132").Value = 2660.3
$ws.Range("K132").Value = 3924.3156
$ws.Range("L132").Value = 7980.900000000001
$ws.Range("M132").Value = -1394.3156
$ws.Range("N132").Value = -13040.9
